$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: search for $needle starting at $cursor.Position (an integer char
# offset) but not beyond $limit, bold it, and return the offset just past
# the match so the next search can continue from there.
# ---------------------------------------------------------------------------
function Bold-Next([int]$start, [int]$limit, [string]$needle) {
    $rng = $d.Range($start, $limit)
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find '$needle' between $start and $limit"
    }
    $rng.Bold = 1
    return $rng.End
}

# ===========================================================================
# 1) Drop the stray _GoBack bookmark that used to sit after the "not a full
#    translation" paragraph.
# ===========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ===========================================================================
# 2) MUST | PHẢI paragraph — bold the six label words.
# ===========================================================================
$p = $d.Paragraphs.Item(5).Range
$pos = $p.Start
$pos = Bold-Next $pos $p.End "MUST"
$pos = Bold-Next $pos $p.End "PHẢI"
$pos = Bold-Next $pos $p.End "REQUIRED"
$pos = Bold-Next $pos $p.End "BẮT-BUỘC"
$pos = Bold-Next $pos $p.End "SHALL"
$pos = Bold-Next $pos $p.End "SẼ"

# ===========================================================================
# 3) MUST NOT paragraph — bold the four label words/phrases.
# ===========================================================================
$p = $d.Paragraphs.Item(6).Range
$pos = $p.Start
$pos = Bold-Next $pos $p.End "MUST NOT"
$pos = Bold-Next $pos $p.End "KHÔNG ĐƯỢC"
$pos = Bold-Next $pos $p.End "SHALL NOT"
$pos = Bold-Next $pos $p.End "SẼ-KHÔNG"

# ===========================================================================
# 4) SHOULD paragraph — bold the label words, then insert " | KHUYẾN-NGHỊ"
#    (bold) just before the closing curly quote.
# ===========================================================================
$p = $d.Paragraphs.Item(7).Range
$pos = $p.Start
$pos = Bold-Next $pos $p.End "SHOULD"
$pos = Bold-Next $pos $p.End "NÊN"
$pos = Bold-Next $pos $p.End "RECOMMENDED"
$pos = Bold-Next $pos $p.End "| KHUYẾN-DÙNG"
$pos = Bold-Next $pos $p.End "KHUYÊN-DÙNG"

# Insert the new " | KHUYẾN-NGHỊ" right after "KHUYÊN-DÙNG" (i.e. right
# before the closing “”) — re-find the closing quote that follows.
$insPoint = $d.Range($pos, $pos)
$insPoint.InsertAfter(" | ")
$pos = $insPoint.End
$newWord = $d.Range($pos, $pos)
$newWord.InsertAfter("KHUYẾN-NGHỊ")
$newWord.Bold = 1

# ===========================================================================
# 5) SHOULD NOT paragraph — bold the label words/phrases (mind the quirky
#    split around "NOT RECOMMENDED | KHÔNG-KHUYẾN-KHÍCH").
# ===========================================================================
$p = $d.Paragraphs.Item(8).Range
$pos = $p.Start
$pos = Bold-Next $pos $p.End "SHOULD NOT"
$pos = Bold-Next $pos $p.End "KHÔNG-NÊN"
$pos = Bold-Next $pos $p.End "NOT RECOMMENDED"
$pos = Bold-Next $pos $p.End " KHÔNG-KHUYẾN-KHÍCH"

# ===========================================================================
# 6) MAY paragraph — bold the label words, and re-wrap "TÙY-CHỌN" with the
#    _GoBack bookmark.
# ===========================================================================
$p = $d.Paragraphs.Item(9).Range
$pos = $p.Start
$pos = Bold-Next $pos $p.End "MAY"
$pos = Bold-Next $pos $p.End "CÓ-THỂ"
$pos = Bold-Next $pos $p.End "OPTIONAL"
$pos = Bold-Next $pos $p.End "TÙY"

$rng = $d.Range($pos, $p.End)
$rng.Find.Execute("TÙY-CHỌN", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$d.Bookmarks.Add("_GoBack", $rng)

Write-Host "done"
